$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "grid on photo" (row 8) as completed
$ws.Range("C8").Value = "y"

# Add new item row for the sign up bug fix
$ws.Range("A17").Value = "signup bug"

# Match the saved selection state from the diff
$ws.Range("B17").Select()
